$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(61221, "Heloísa Nogueira", "Atendimento ao Cliente", "Doença", 3, 45103, 9654.08),
    @(69785, "Clara Silveira", "Marketing", "Outros", 7, 45100, 2519.69),
    @(79393, "Kaique da Rosa", "Recursos Humanos", "Problemas pessoais", 2, 45099, 5223.71),
    @(61715, "Nicolas Alves", "Recursos Humanos", "Doença", 1, 45097, 11926.95),
    @(27968, "Stella Moura", "TI", "Doença", 6, 45102, 5208.97),
    @(28414, "Ana Laura Ferreira", "Jurídico", "Consulta médica", 5, 45088, 3992.54),
    @(7685, "Nathan Costa", "Marketing", "Outros", 3, 45090, 9238.93),
    @(38304, "Marcela Rocha", "Recursos Humanos", "Viagem de negócios", 6, 45096, 7841.25),
    @(30070, "Emanuella Costa", "Vendas", "Consulta médica", 7, 45081, 9840.77),
    @(33795, "Alice Pires", "Jurídico", "Doença", 2, 45084, 9759.25)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
